# Update column F (dSF) values for specific rows as part of a data repull / mean calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -2
    17 = 2
    21 = 2
    26 = -2
    33 = 4
    39 = 3
    41 = -1
    42 = -1
    44 = -1
    46 = 0
    51 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
